# Insert a new weekly price record as row 95 (Mango, "Especial" quality,
# Agricola del Norte S.A. de Arica), shifting the existing rows 95-118 down
# to 96-119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("95:95").Insert()

$ws.Cells.Item(95, 1).Value  = 1
$ws.Cells.Item(95, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(95, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(95, 4).Value  = 44572
$ws.Cells.Item(95, 5).Value  = 15
$ws.Cells.Item(95, 6).Value  = "Fruta"
$ws.Cells.Item(95, 7).Value  = 100108
$ws.Cells.Item(95, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(95, 9).Value  = 100108002
$ws.Cells.Item(95, 10).Value = "Mango"
$ws.Cells.Item(95, 11).Value = "Sin especificar"
$ws.Cells.Item(95, 12).Value = "Especial"
$ws.Cells.Item(95, 13).Value = 456
$ws.Cells.Item(95, 14).Value = 6000
$ws.Cells.Item(95, 15).Value = 6500
$ws.Cells.Item(95, 16).Value = 6250
$ws.Cells.Item(95, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(95, 18).Value = "Perú"
$ws.Cells.Item(95, 19).Value = 1562
$ws.Cells.Item(95, 20).Value = 4
